$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "D2" = "246.13"
    "D3" = "24.19"
    "D4" = "5.287"
    "D6" = "6.500"
    "D7" = "3.143"
    "D8" = "0.8121"
    "D9" = "0.8592"
    "D11" = "0.06993"
    "D12" = "0.03135"
    "D13" = "0.02933"
    "D14" = "0.09402"
    "D15" = "3.773"
    "D16" = "0.001526"
    "D17" = "0.04677"
    "B18" = "TigerCash"
    "C18" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "D18" = "0.006091"
    "E18" = "17TigerCashTCH"
    "B19" = "BitKan"
    "C19" = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
    "D19" = "0.001238"
    "E19" = "18BitKanKAN"
    "B20" = "HotbitToken"
    "C20" = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
    "D20" = "0.004640"
    "E20" = "19HotbitTokenHTB"
    "B21" = "NitroEx"
    "C21" = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
    "D21" = "0.00006102"
    "E21" = "20NitroExNTXWorstin24h"
    "B22" = "LEO"
    "C22" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "D22" = "3.503"
    "E22" = "21LEOLEO"
    "B23" = "BTSEToken"
    "C23" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
    "D23" = "2.137"
    "E23" = "22BTSETokenBTSE"
    "B24" = "One"
    "C24" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "D24" = "0.009862"
    "E24" = "23OneONEBestin24h"
    "D28" = "0.0002332"
    "E28" = "27UpBotsUBXT"
    "D40" = "0.03704"
    "B41" = "KickToken"
    "C41" = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
    "D41" = "0.006262"
    "E41" = "40KickTokenKICK"
    "B42" = "BKEXToken"
    "C42" = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
    "D42" = "0.1058"
    "E42" = "41BKEXTokenBKK"
    "B43" = "CEJI"
    "C43" = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
    "D43" = "0.003401"
    "E43" = "42CEJICEJI"
    "D44" = "0.008516"
    "D45" = "0.00005283"
    "D47" = "0.4401"
    "D48" = "0.002404"
    "E48" = "47BOLOBOLO"
    "D49" = "0.00002101"
    "D50" = "0.0002001"
}

foreach ($ref in $changes.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$ref]
    $cell.Style = "Normal"
}
